# PCB Technology Requirement Format -> ULTiM8x8 bus board revision
# (UlTiM_bus layout, gerber complete. Bus accessory for ULTiM8x8)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PCB Name & Version
$ws.Range("B3").Value = "ULTiM8x8"

# PCB Size (inches)
$ws.Range("C5").Value = "( 2.40 ) × (2.40 )"

# Thickness
$ws.Range("B8").Value = "1.6mm"

# Copper Weight: value moves from C9 to B9, and changes 2oz -> 1oz
$ws.Range("C9").Value = ""
$ws.Range("B9").Value = "1oz"

# PCB Colour: value moves from D11 to C10, "other:( CLEAR )" -> "White白色"
$ws.Range("D11").Value = ""
$ws.Range("C10").Value = "White白色"

# Minimum trace/space & Minimum hole size values
$ws.Range("B15").Value = "0.254mm/0.254mm"
$ws.Range("D15").Value = "0.635mm"
$ws.Range("D15").HorizontalAlignment = -4108

# Special requirements: " Yes" -> "No"
$ws.Range("B16").Value = "No"

# Comment: remove the silkscreen note
$ws.Range("B17").Value = ""
